# Update countries & provincias Spain
# Applies the daily data refresh to the "Pais" sheet:
#  - bumps the "Datos actualizados" timestamp
#  - updates case counters for several countries
#  - Hungria/Albania swap places (Hungria now gets the refreshed numbers,
#    Albania keeps the figures that used to belong to the other row)
#  - Georgia/Yemen swap places the same way

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 09:12"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6636266
$ws.Range("C4").Value = 19
$ws.Range("E4").Value = 2520883

# Armenia (row 62)
$ws.Range("B62").Value = 45675
$ws.Range("C62").Value = 172
$ws.Range("D62").Value = 41605
$ws.Range("E62").Value = 3159
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 911

# El Salvador (row 74)
$ws.Range("B74").Value = 26851
$ws.Range("C74").Value = 78
$ws.Range("D74").Value = 17874
$ws.Range("E74").Value = 8195

# Row 93 becomes Hungria with refreshed figures
$ws.Range("A93").Value = "Hungria"
$ws.Range("B93").Value = 11825
$ws.Range("C93").Value = 916
$ws.Range("D93").Value = 4058
$ws.Range("E93").Value = 7134
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 633

# Row 94 becomes Albania, carrying the figures previously on row 93
$ws.Range("A94").Value = "Albania"
$ws.Range("B94").Value = 11021
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 6443
$ws.Range("E94").Value = 4251
$ws.Range("H94").Value = 327

# Birmania (row 144)
$ws.Range("B144").Value = 2445
$ws.Range("C144").Value = 23
$ws.Range("E144").Value = 1806

# Row 151 becomes Georgia with refreshed figures
$ws.Range("A151").Value = "Georgia"
$ws.Range("B151").Value = 2075
$ws.Range("C151").Value = 158
$ws.Range("D151").Value = 1363
$ws.Range("E151").Value = 693
$ws.Range("H151").Value = 19

# Row 152 becomes Yemen, carrying the figures previously on row 151
$ws.Range("A152").Value = "Yemen"
$ws.Range("B152").Value = 2007
$ws.Range("D152").Value = 1211
$ws.Range("E152").Value = 214
$ws.Range("H152").Value = 582

# Letonia (row 159)
$ws.Range("B159").Value = 1464
$ws.Range("C159").Value = 5
$ws.Range("E159").Value = 181
